# Hortaliza, Vega Modelo de Temuco - Achicoria: add two new weekly price rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new row at sheet row 18 (pushes the old rows 18-26 down by one) ---
$ws.Rows("18:18").Insert()

$ws.Cells.Item(18, 1).Value = 10
$ws.Cells.Item(18, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(18, 3).Value = "La Araucanía"
$ws.Cells.Item(18, 4).Value = 44679
$ws.Cells.Item(18, 5).Value = 9
$ws.Cells.Item(18, 6).Value = 100112010
$ws.Cells.Item(18, 7).Value = "Achicoria"
$ws.Cells.Item(18, 8).Value = "Sin especificar"
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 90
$ws.Cells.Item(18, 11).Value = 12000
$ws.Cells.Item(18, 12).Value = 12000
$ws.Cells.Item(18, 13).Value = 12000
$ws.Cells.Item(18, 14).Value = "`$/caja 18 unidades"
$ws.Cells.Item(18, 15).Value = "Región Metropolitana"
$ws.Cells.Item(18, 16).Value = 667
$ws.Cells.Item(18, 17).Value = 18
$ws.Cells.Item(18, 18).Value = "Hortaliza"

# --- Insert another new row at sheet row 26 (after the row that used to be 25) ---
$ws.Rows("26:26").Insert()

$ws.Cells.Item(26, 1).Value = 10
$ws.Cells.Item(26, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(26, 3).Value = "La Araucanía"
$ws.Cells.Item(26, 4).Value = 44301
$ws.Cells.Item(26, 5).Value = 9
$ws.Cells.Item(26, 6).Value = 100112010
$ws.Cells.Item(26, 7).Value = "Achicoria"
$ws.Cells.Item(26, 8).Value = "Sin especificar"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 40
$ws.Cells.Item(26, 11).Value = 12000
$ws.Cells.Item(26, 12).Value = 12000
$ws.Cells.Item(26, 13).Value = 12000
$ws.Cells.Item(26, 14).Value = "`$/caja 16 unidades"
$ws.Cells.Item(26, 15).Value = "Región del Maule"
$ws.Cells.Item(26, 16).Value = 750
$ws.Cells.Item(26, 17).Value = 16
$ws.Cells.Item(26, 18).Value = "Hortaliza"

$ws.Range("A1").Select()
